$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C ("Förändrad") from 45735 to 45736 for all data rows (2-44)
for ($r = 2; $r -le 44; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45735) {
        $cell.Value2 = 45736
    }
}
